$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the panel-query timestamps on the "data" sheet (F2:F5) ---
$ws.Range("F2").Value() = "2021-10-05 14:19:18.637625"
$ws.Range("F3").Value() = "2021-10-05 14:19:18.637634"
$ws.Range("F4").Value() = "2021-10-05 14:19:18.637637"
$ws.Range("F5").Value() = "2021-10-05 14:19:18.637640"

# --- Add a new "metadata" sheet right after "data" ---
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$newSheet.Name = "metadata"

# Match the page margins used by the rest of the workbook.
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Header row (B1:G1) - copy formatting (bold/border/centered) from the
# "data" sheet's header row, then overwrite the text for each column.
$ws.Range("B1:F1").Copy($newSheet.Range("B1"))
$ws.Range("F1").Copy($newSheet.Range("G1"))

$newSheet.Range("B1").Value() = "data_name"
$newSheet.Range("C1").Value() = "data_id"
$newSheet.Range("D1").Value() = "data_version"
$newSheet.Range("E1").Value() = "data_version_created"
$newSheet.Range("F1").Value() = "panel_query_time"
$newSheet.Range("G1").Value() = "panel_get_request"

# Data row (row 2). A2 keeps the same index styling as on "data".
$ws.Range("A2").Copy($newSheet.Range("A2"))
$newSheet.Range("A2").Value() = 0

$newSheet.Range("B2").Value() = "Bladder cancer pertinent cancer susceptibility"
$newSheet.Range("C2").Value() = 208

# Force "1.1" to be stored as text (not coerced to a number), then drop
# back to the default (unstyled) cell format to match the other cells.
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value() = "1.1"
$newSheet.Range("D2").Style = "Normal"

$newSheet.Range("E2").Value() = "2019-06-20T15:10:18.540382Z"
$newSheet.Range("F2").Value() = "2021-10-05 14:19:18.633680"
$newSheet.Range("G2").Value() = "https://panelapp.genomicsengland.co.uk/api/v1/panels/208/?format=json"

# Keep "data" as the active/selected sheet (unchanged by the diff).
$ws.Activate() | Out-Null
$ws.Range("A1").Select() | Out-Null
